# Feat: Player HoverState 추가 및 State 로직 수정
# Adds a new "hovering" field/column to the PartDbSheet table, between
# "smoothRotation" and "bulletPrefab_Path". All columns from the old
# "bulletPrefab_Path" (N) through "isSplash" (W) shift one column to the
# right (O..X) to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before N; this shifts N:W -> O:X and keeps every
# existing cell's style/type intact (matches Excel's native Insert Column
# behavior, including the width the new column inherits from its left
# neighbor).
$ws.Columns("N:N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# New header for the inserted column.
$ws.Range("N1").Value = "hovering"

# New per-part "hovering" values (numeric), one per data row.
$ws.Range("N2").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("N4").Value = 0.1
$ws.Range("N5").Value = 0.1
$ws.Range("N6").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("N12").Value = 0

# Restore the view/selection state recorded for this edit.
$ws.Range("N6").Select()
